# Applies the data corrections described in the commit diff for the
# "Jogos_da_Semana_FlashScore_2024-10-31.xlsx" workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2 updates ---
$ws.Range("G2").Value  = 2.35
$ws.Range("I2").Value  = 3.25
$ws.Range("M2").Value  = 1.1
$ws.Range("N2").Value  = 7
$ws.Range("X2").Value  = 10
$ws.Range("Y2").Value  = 10
$ws.Range("AB2").Value = 41
$ws.Range("AJ2").Value = 34
$ws.Range("AK2").Value = 29
$ws.Range("AN2").Value = 4.33
$ws.Range("AP2").Value = 29
$ws.Range("AQ2").Value = 51
$ws.Range("AX2").Value = 19

# --- Row 3 updates ---
$ws.Range("BD3").Value = 126
